{"js": "/*\n * Weekly math worksheet refresh:\n *   - Update the date heading to the new date.\n *   - Replace the 100 addition/subtraction problems in the 20x5 table\n *     (one problem per cell, in row-major reading order) with the new\n *     set of problems.\n *\n * Each cell/paragraph keeps its original formatting (fonts, size,\n * alignment) because we replace the *text of a Range* rather than\n * replacing the whole cell/paragraph body, which would drop the\n * existing run/paragraph properties.\n */\n\n// 1) Update the date heading (first paragraph in the document body).\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst dateParagraph = paragraphs.items[0];\ndateParagraph\n  .getRange(\"Whole\")\n  .insertText(\"2025-03-01 Saturday\", Word.InsertLocation.replace);\nawait context.sync();\n\n// 2) Update every arithmetic problem in the single table, row by row,\n//    left to right - matching the order the problems appear in the file.\nconst newValues = [\n  [\"63-7=\", \"89-45=\", \"74+6=\", \"18+9=\", \"77-30=\"],\n  [\"56+0=\", \"36+56=\", \"89-50=\", \"29+39=\", \"11+86=\"],\n  [\"73-40=\", \"76+13=\", \"79-66=\", \"73-69=\", \"48+3=\"],\n  [\"98-51=\", \"56-0=\", \"72+0=\", \"20+4=\", \"36+50=\"],\n  [\"91-83=\", \"48-25=\", \"46+13=\", \"52-34=\", \"85-16=\"],\n  [\"10+32=\", \"36-19=\", \"95-84=\", \"56+33=\", \"22+22=\"],\n  [\"7+4=\", \"49+29=\", \"49-35=\", \"54-27=\", \"1+39=\"],\n  [\"37+16=\", \"70-49=\", \"70-9=\", \"68-55=\", \"93-75=\"],\n  [\"1+23=\", \"7+11=\", \"12+14=\", \"8+58=\", \"51+15=\"],\n  [\"75-5=\", \"1+1=\", \"87-79=\", \"35+45=\", \"69-16=\"],\n  [\"18+57=\", \"49+32=\", \"23+43=\", \"57-12=\", \"27+25=\"],\n  [\"99-28=\", \"93-4=\", \"56+14=\", \"15+27=\", \"31-30=\"],\n  [\"66-35=\", \"93+0=\", \"97-48=\", \"39-19=\", \"20+50=\"],\n  [\"57+3=\", \"88-54=\", \"46+4=\", \"2+77=\", \"17+77=\"],\n  [\"90-28=\", \"45-39=\", \"83-13=\", \"34-18=\", \"24+58=\"],\n  [\"68-54=\", \"99-91=\", \"16+70=\", \"15+78=\", \"95-59=\"],\n  [\"20+50=\", \"42-41=\", \"72-20=\", \"35-8=\", \"32-21=\"],\n  [\"39+43=\", \"36+54=\", \"76-25=\", \"91-44=\", \"75-7=\"],\n  [\"8+51=\", \"96-8=\", \"25+31=\", \"90+4=\", \"18-10=\"],\n  [\"26+45=\", \"95-8=\", \"51-21=\", \"64-2=\", \"73-37=\"]\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.load(\"rowCount,values\");\nawait context.sync();\n\nconst columnCount = table.values.length > 0 ? table.values[0].length : 0;\n\nfor (let row = 0; row < table.rowCount; row++) {\n  for (let col = 0; col < columnCount; col++) {\n    const cell = table.getCell(row, col);\n    cell.body\n      .getRange(\"Whole\")\n      .insertText(newValues[row][col], Word.InsertLocation.replace);\n  }\n}\nawait context.sync();\n", "ps1": "# Weekly math worksheet refresh:\n#   - Update the date heading to the new date.\n#   - Replace the 100 addition/subtraction problems in the 20x5 table\n#     (one problem per cell, in row-major reading order) with the new\n#     set of problems.\n#\n# Assigning directly to a Range's .Text preserves the existing run /\n# paragraph formatting (fonts, size, alignment) of that range, because\n# only the text content is replaced - not the run or paragraph objects.\n\n$d = $word.ActiveDocument\n\n# 1) Update the date heading (first paragraph in the document body).\n$d.Paragraphs.Item(1).Range.Text = \"2025-03-01 Saturday\"\n\n# 2) Update every arithmetic problem in the single table, row by row,\n#    left to right - matching the order the problems appear in the file.\n$newValues = @(\n    @(\"63-7=\", \"89-45=\", \"74+6=\", \"18+9=\", \"77-30=\"),\n    @(\"56+0=\", \"36+56=\", \"89-50=\", \"29+39=\", \"11+86=\"),\n    @(\"73-40=\", \"76+13=\", \"79-66=\", \"73-69=\", \"48+3=\"),\n    @(\"98-51=\", \"56-0=\", \"72+0=\", \"20+4=\", \"36+50=\"),\n    @(\"91-83=\", \"48-25=\", \"46+13=\", \"52-34=\", \"85-16=\"),\n    @(\"10+32=\", \"36-19=\", \"95-84=\", \"56+33=\", \"22+22=\"),\n    @(\"7+4=\", \"49+29=\", \"49-35=\", \"54-27=\", \"1+39=\"),\n    @(\"37+16=\", \"70-49=\", \"70-9=\", \"68-55=\", \"93-75=\"),\n    @(\"1+23=\", \"7+11=\", \"12+14=\", \"8+58=\", \"51+15=\"),\n    @(\"75-5=\", \"1+1=\", \"87-79=\", \"35+45=\", \"69-16=\"),\n    @(\"18+57=\", \"49+32=\", \"23+43=\", \"57-12=\", \"27+25=\"),\n    @(\"99-28=\", \"93-4=\", \"56+14=\", \"15+27=\", \"31-30=\"),\n    @(\"66-35=\", \"93+0=\", \"97-48=\", \"39-19=\", \"20+50=\"),\n    @(\"57+3=\", \"88-54=\", \"46+4=\", \"2+77=\", \"17+77=\"),\n    @(\"90-28=\", \"45-39=\", \"83-13=\", \"34-18=\", \"24+58=\"),\n    @(\"68-54=\", \"99-91=\", \"16+70=\", \"15+78=\", \"95-59=\"),\n    @(\"20+50=\", \"42-41=\", \"72-20=\", \"35-8=\", \"32-21=\"),\n    @(\"39+43=\", \"36+54=\", \"76-25=\", \"91-44=\", \"75-7=\"),\n    @(\"8+51=\", \"96-8=\", \"25+31=\", \"90+4=\", \"18-10=\"),\n    @(\"26+45=\", \"95-8=\", \"51-21=\", \"64-2=\", \"73-37=\")\n)\n\n$table = $d.Tables.Item(1)\n$rowCount = $table.Rows.Count\n$colCount = $table.Columns.Count\n\nfor ($r = 1; $r -le $rowCount; $r++) {\n    for ($c = 1; $c -le $colCount; $c++) {\n        $table.Cell($r, $c).Range.Text = $newValues[$r - 1][$c - 1]\n    }\n}\n"}
